# Update mods data [2025-11-22 15:08:32]
#
# Appends a new daily data row (row 13) to the ModCounts sheet:
#   Date = 2025/11/22, Game = 逃离鸭科夫, ModCount = 1222
# matching the style/format already used by the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 13
$prevRow = $newRow - 1

# Column A holds a date-like string ("2025/11/11", "2025/11/12", ...) that is
# stored as plain literal text in every existing row, NOT as a real Excel
# date. Force a text number format before assigning the value so Excel does
# not auto-parse "2025/11/22" into a date serial number.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2025/11/22"
$ws.Range("B$newRow").Value = "逃离鸭科夫"
$ws.Range("C$newRow").Value = 1222

# Match the formatting (centered alignment, etc.) used by the rest of the
# data rows by copying the previous row's format onto the new row.
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$ws.Range("A${prevRow}:C$prevRow").Copy()
$ws.Range("A${newRow}:C$newRow").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
